# Update the TEMPLATE 2026 GWFMOA docx:
#  1. Bump the empty "Recent News" list-paragraph's left indent from 76 -> 796 twips.
#  2. Insert a new bullet abstractNum (nsid 24354785) at abstractNumId="1",
#     renumbering the existing abstractNum ids 1..4 -> 2..5, and add a new
#     <w:num numId="6"> that points at the freshly inserted abstractNum.
#
# Paragraph-format property setters (ParagraphFormat.LeftIndent, etc.) are not
# reliable on paragraphs that live inside a table cell in this host, and the
# numbering part has no dedicated object-model surface for raw abstractNum
# insertion/renumbering, so this is done as a precise text edit of the
# document's flattened WordOpenXML (round-tripped back through the same
# property), which is the supported mechanism for this kind of bulk/structural
# part edit.

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

# --- 1. document.xml: ind left 76 -> 796 on the empty ListParagraph in the
#        "Recent News" table cell.
$oldInd = 'w:ind w:left="76"/'
$newInd = 'w:ind w:left="796"/'
if ($xml.IndexOf($oldInd) -lt 0) {
    throw "expected indent anchor not found"
}
$xml = $xml.Replace($oldInd, $newInd)

# --- 2. word/numbering.xml: insert the new abstractNum + renumber existing ones.

# New abstractNum body (nsid 24354785), to be inserted as abstractNumId="1".
$newAbstractNumBody = '<w:nsid w:val="24354785"/><w:multiLevelType w:val="hybridMultilevel"/><w:tmpl w:val="D5A820FA"/><w:lvl w:ilvl="0" w:tplc="04090001"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="'+[char]0xF0B7+'"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="796" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="1" w:tplc="04090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="1516" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="2" w:tplc="04090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="'+[char]0xF0A7+'"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2236" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="3" w:tplc="04090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="'+[char]0xF0B7+'"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="2956" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="4" w:tplc="04090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="3676" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="5" w:tplc="04090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="'+[char]0xF0A7+'"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="4396" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="6" w:tplc="04090001" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="'+[char]0xF0B7+'"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5116" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Symbol" w:hAnsi="Symbol" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="7" w:tplc="04090003" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="o"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="5836" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New" w:hint="default"/></w:rPr></w:lvl><w:lvl w:ilvl="8" w:tplc="04090005" w:tentative="1"><w:start w:val="1"/><w:numFmt w:val="bullet"/><w:lvlText w:val="'+[char]0xF0A7+'"/><w:lvlJc w:val="left"/><w:pPr><w:ind w:left="6556" w:hanging="360"/></w:pPr><w:rPr><w:rFonts w:ascii="Wingdings" w:hAnsi="Wingdings" w:hint="default"/></w:rPr></w:lvl>'

# abstractNumId="1" (nsid 28CD0207) becomes abstractNumId="2" ; the new
# abstractNum body above is inserted in the now-vacated id="1" slot, right
# before it.
$anchor1 = '<w:abstractNum w:abstractNumId="1" w15:restartNumberingAfterBreak="0"><w:nsid w:val="28CD0207"'
$replacement1 = '<w:abstractNum w:abstractNumId="1" w15:restartNumberingAfterBreak="0">' + $newAbstractNumBody + '</w:abstractNum><w:abstractNum w:abstractNumId="2" w15:restartNumberingAfterBreak="0"><w:nsid w:val="28CD0207"'
if ($xml.IndexOf($anchor1) -lt 0) {
    throw "abstractNum anchor 28CD0207 not found"
}
$xml = $xml.Replace($anchor1, $replacement1)

# abstractNumId="2" (nsid 47F30842) becomes abstractNumId="3"
$anchor2 = '<w:abstractNum w:abstractNumId="2" w15:restartNumberingAfterBreak="0"><w:nsid w:val="47F30842"'
$replacement2 = '<w:abstractNum w:abstractNumId="3" w15:restartNumberingAfterBreak="0"><w:nsid w:val="47F30842"'
if ($xml.IndexOf($anchor2) -lt 0) {
    throw "abstractNum anchor 47F30842 not found"
}
$xml = $xml.Replace($anchor2, $replacement2)

# abstractNumId="3" (nsid 49896415) becomes abstractNumId="4"
$anchor3 = '<w:abstractNum w:abstractNumId="3" w15:restartNumberingAfterBreak="0"><w:nsid w:val="49896415"'
$replacement3 = '<w:abstractNum w:abstractNumId="4" w15:restartNumberingAfterBreak="0"><w:nsid w:val="49896415"'
if ($xml.IndexOf($anchor3) -lt 0) {
    throw "abstractNum anchor 49896415 not found"
}
$xml = $xml.Replace($anchor3, $replacement3)

# abstractNumId="4" (nsid 53D8540D) becomes abstractNumId="5"
$anchor4 = '<w:abstractNum w:abstractNumId="4" w15:restartNumberingAfterBreak="0"><w:nsid w:val="53D8540D"'
$replacement4 = '<w:abstractNum w:abstractNumId="5" w15:restartNumberingAfterBreak="0"><w:nsid w:val="53D8540D"'
if ($xml.IndexOf($anchor4) -lt 0) {
    throw "abstractNum anchor 53D8540D not found"
}
$xml = $xml.Replace($anchor4, $replacement4)

# --- <w:num> map: keep numId 1..5 pointing at the same logical lists (now
#     shifted abstractNumIds) and add a new numId="6" for the inserted list.
$oldNums = '<w:num w:numId="1" w16cid:durableId="1273973645"><w:abstractNumId w:val="0"/></w:num><w:num w:numId="2" w16cid:durableId="937832032"><w:abstractNumId w:val="3"/></w:num><w:num w:numId="3" w16cid:durableId="147747548"><w:abstractNumId w:val="1"/></w:num><w:num w:numId="4" w16cid:durableId="1685084077"><w:abstractNumId w:val="2"/></w:num><w:num w:numId="5" w16cid:durableId="1129973996"><w:abstractNumId w:val="4"/></w:num>'
$newNums = '<w:num w:numId="1" w16cid:durableId="1273973645"><w:abstractNumId w:val="0"/></w:num><w:num w:numId="2" w16cid:durableId="937832032"><w:abstractNumId w:val="4"/></w:num><w:num w:numId="3" w16cid:durableId="147747548"><w:abstractNumId w:val="2"/></w:num><w:num w:numId="4" w16cid:durableId="1685084077"><w:abstractNumId w:val="3"/></w:num><w:num w:numId="5" w16cid:durableId="1129973996"><w:abstractNumId w:val="5"/></w:num><w:num w:numId="6" w16cid:durableId="184369170"><w:abstractNumId w:val="1"/></w:num>'
if ($xml.IndexOf($oldNums) -lt 0) {
    throw "num map anchor not found"
}
$xml = $xml.Replace($oldNums, $newNums)

$d.WordOpenXML = $xml
Write-Output "ok"
